$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D2"  = 10.1403762618524
    "E2"  = 0.3401649376635241
    "D3"  = 10.0144558299543
    "E3"  = 0.3295199693804857
    "D4"  = 10.97473529764362
    "E4"  = 0.3392952363579432
    "D5"  = 11.27740638296911
    "E5"  = 0.352369311859356
    "D6"  = 11.61404378698954
    "E6"  = 0.3363480743703146
    "D7"  = 12.42598771940949
    "E7"  = 0.3776266911159403
    "D8"  = 12.25033722608613
    "E8"  = 0.3298622200525696
    "D9"  = 13.28187451693018
    "E9"  = 0.4918318573854879
    "D10" = 13.04722761421902
    "E10" = 0.3532769108607338
    "D11" = 14.3623631510613
    "E11" = 0.4785707621321235
    "D12" = 14.02456987818324
    "E12" = 0.3981500633966313
    "D13" = 15.31754357659524
    "E13" = 0.4992000233112369
    "D14" = 14.79142891025522
    "E14" = 0.4317532294218546
    "D15" = 15.96220731988605
    "E15" = 0.4902791803780185
    "D16" = 15.49040288444248
    "E16" = 0.4565458398706101
    "D17" = 16.52404572044733
    "E17" = 0.4898259902740931
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
